# Insert a new data row at row 28 (weekly price update), pushing all
# subsequent rows (28-102) down by one (to 29-103).  The new row 28
# contains a fresh "Madrigal" / "Primera" quote for 2022-08-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("28:28").Insert()

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44804
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 100112013
$ws.Range("G28").Value = "Alcachofa"
$ws.Range("H28").Value = "Madrigal"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 13000
$ws.Range("L28").Value = 13000
$ws.Range("M28").Value = 13000
$ws.Range("N28").Value = "$/caja 40 unidades"
$ws.Range("O28").Value = "Provincia del Elquí"
$ws.Range("P28").Value = 325
$ws.Range("Q28").Value = 40
$ws.Range("R28").Value = "Hortaliza"
